$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1143380298"
$ws.Range("D16").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E16").Value = "2012"
$ws.Range("F16").Value = 230066
$ws.Range("G16").Value = 14500000
$ws.Range("C17").Value = "1143380298"
$ws.Range("D17").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E17").Value = "2101"
$ws.Range("F17").Value = 406000
$ws.Range("G17").Value = 14500000
$ws.Range("C18").Value = "1143380298"
$ws.Range("D18").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E18").Value = "2102"
$ws.Range("F18").Value = 406000
$ws.Range("G18").Value = 14500000
$ws.Range("C19").Value = "1143380298"
$ws.Range("D19").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E19").Value = "2103"
$ws.Range("F19").Value = 406000
$ws.Range("G19").Value = 14500000
$ws.Range("C20").Value = "1143362651"
$ws.Range("D20").Value = "CARLOS HUMBERTO VILLA MERCADO"
$ws.Range("E20").Value = "2104"
$ws.Range("F20").Value = 580000
$ws.Range("G20").Value = 17127677
$ws.Range("C21").Value = "1143380298"
$ws.Range("D21").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E21").Value = "2104"
$ws.Range("F21").Value = 406000
$ws.Range("G21").Value = 14500000
$ws.Range("C22").Value = "1143362651"
$ws.Range("D22").Value = "CARLOS HUMBERTO VILLA MERCADO"
$ws.Range("E22").Value = "2105"
$ws.Range("F22").Value = 580000
$ws.Range("G22").Value = 17127677
$ws.Range("C23").Value = "1143380298"
$ws.Range("D23").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E23").Value = "2105"
$ws.Range("F23").Value = 406000
$ws.Range("G23").Value = 14500000
$ws.Range("C24").Value = "1143362651"
$ws.Range("D24").Value = "CARLOS HUMBERTO VILLA MERCADO"
$ws.Range("E24").Value = "2106"
$ws.Range("F24").Value = 580000
$ws.Range("G24").Value = 17127677
$ws.Range("C25").Value = "1143380298"
$ws.Range("D25").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E25").Value = "2106"
$ws.Range("F25").Value = 406000
$ws.Range("G25").Value = 14500000
$ws.Range("C26").Value = "1143362651"
$ws.Range("D26").Value = "CARLOS HUMBERTO VILLA MERCADO"
$ws.Range("E26").Value = "2107"
$ws.Range("F26").Value = 580000
$ws.Range("G26").Value = 17127677
$ws.Range("C27").Value = "1143380298"
$ws.Range("D27").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E27").Value = "2107"
$ws.Range("F27").Value = 406000
$ws.Range("G27").Value = 14500000
$ws.Range("C28").Value = "1143380298"
$ws.Range("D28").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E28").Value = "2108"
$ws.Range("F28").Value = 406000
$ws.Range("G28").Value = 14500000
$ws.Range("C29").Value = "1143380298"
$ws.Range("D29").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E29").Value = "2109"
$ws.Range("F29").Value = 406000
$ws.Range("G29").Value = 14500000
$ws.Range("C30").Value = "1143380298"
$ws.Range("D30").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E30").Value = "2110"
$ws.Range("F30").Value = 406000
$ws.Range("G30").Value = 14500000
$ws.Range("C31").Value = "1143380298"
$ws.Range("D31").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E31").Value = "2111"
$ws.Range("F31").Value = 406000
$ws.Range("G31").Value = 14500000
$ws.Range("C32").Value = "1143380298"
$ws.Range("D32").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E32").Value = "2112"
$ws.Range("F32").Value = 406000
$ws.Range("G32").Value = 14500000
$ws.Range("C33").Value = "1143380298"
$ws.Range("D33").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E33").Value = "2201"
$ws.Range("F33").Value = 406000
$ws.Range("G33").Value = 14500000
$ws.Range("C34").Value = "1143380298"
$ws.Range("D34").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E34").Value = "2202"
$ws.Range("F34").Value = 406000
$ws.Range("G34").Value = 14500000
$ws.Range("C35").Value = "1143380298"
$ws.Range("D35").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E35").Value = "2203"
$ws.Range("F35").Value = 406000
$ws.Range("G35").Value = 14500000
$ws.Range("C36").Value = "1143380298"
$ws.Range("D36").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E36").Value = "2204"
$ws.Range("F36").Value = 406000
$ws.Range("G36").Value = 14500000
$ws.Range("C37").Value = "1143380298"
$ws.Range("D37").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E37").Value = "2205"
$ws.Range("F37").Value = 406000
$ws.Range("G37").Value = 14500000
$ws.Range("C38").Value = "1143380298"
$ws.Range("D38").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E38").Value = "2206"
$ws.Range("F38").Value = 406000
$ws.Range("G38").Value = 14500000
$ws.Range("C39").Value = "1143380298"
$ws.Range("D39").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E39").Value = "2207"
$ws.Range("F39").Value = 406000
$ws.Range("G39").Value = 14500000
$ws.Range("C40").Value = "1143380298"
$ws.Range("D40").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E40").Value = "2208"
$ws.Range("F40").Value = 406000
$ws.Range("G40").Value = 14500000
$ws.Range("C41").Value = "1143380298"
$ws.Range("D41").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E41").Value = "2209"
$ws.Range("F41").Value = 406000
$ws.Range("G41").Value = 14500000
$ws.Range("C42").Value = "1143380298"
$ws.Range("D42").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E42").Value = "2210"
$ws.Range("F42").Value = 406000
$ws.Range("G42").Value = 14500000
$ws.Range("C43").Value = "1143380298"
$ws.Range("D43").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E43").Value = "2211"
$ws.Range("F43").Value = 406000
$ws.Range("G43").Value = 14500000
$ws.Range("C44").Value = "1143380298"
$ws.Range("D44").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E44").Value = "2212"
$ws.Range("F44").Value = 406000
$ws.Range("G44").Value = 14500000
$ws.Range("C45").Value = "1143380298"
$ws.Range("D45").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E45").Value = "2301"
$ws.Range("F45").Value = 406000
$ws.Range("G45").Value = 14500000
$ws.Range("C46").Value = "1143380298"
$ws.Range("D46").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E46").Value = "2302"
$ws.Range("F46").Value = 406000
$ws.Range("G46").Value = 14500000
$ws.Range("C47").Value = "1143380298"
$ws.Range("D47").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E47").Value = "2303"
$ws.Range("F47").Value = 406000
$ws.Range("G47").Value = 14500000
$ws.Range("C48").Value = "1143380298"
$ws.Range("D48").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E48").Value = "2304"
$ws.Range("F48").Value = 406000
$ws.Range("G48").Value = 14500000
$ws.Range("C49").Value = "1143380298"
$ws.Range("D49").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E49").Value = "2305"
$ws.Range("F49").Value = 406000
$ws.Range("G49").Value = 14500000
$ws.Range("C50").Value = "1143380298"
$ws.Range("D50").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E50").Value = "2306"
$ws.Range("F50").Value = 406000
$ws.Range("G50").Value = 14500000
$ws.Range("C51").Value = "1143380298"
$ws.Range("D51").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E51").Value = "2307"
$ws.Range("F51").Value = 406000
$ws.Range("G51").Value = 14500000
$ws.Range("C52").Value = "1143380298"
$ws.Range("D52").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E52").Value = "2308"
$ws.Range("F52").Value = 406000
$ws.Range("G52").Value = 14500000
$ws.Range("C53").Value = "1143380298"
$ws.Range("D53").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E53").Value = "2309"
$ws.Range("F53").Value = 406000
$ws.Range("G53").Value = 14500000
$ws.Range("C54").Value = "1143380298"
$ws.Range("D54").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E54").Value = "2310"
$ws.Range("F54").Value = 406000
$ws.Range("G54").Value = 14500000
$ws.Range("C55").Value = "1143380298"
$ws.Range("D55").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E55").Value = "2311"
$ws.Range("F55").Value = 406000
$ws.Range("G55").Value = 14500000
$ws.Range("C56").Value = "1143380298"
$ws.Range("D56").Value = "ESTEBAN DARIO BARBOZA PRIMERA"
$ws.Range("E56").Value = "2312"
$ws.Range("F56").Value = 257133
$ws.Range("G56").Value = 14500000
